# Reroute test data update
# QA_538.xlsx - "Input" sheet carries OrderId / CloneOrderId / RecurringExpiry
# sample values used by the CustomOrder staging test. The old order/clone ids
# (and the recurring-expiry date that went with them) have been retired, so
# swap them for the freshly issued ids/date.
#
# These columns hold digit-only / date-shaped strings but must stay TEXT
# (they are opaque identifiers, not numbers to do arithmetic on), so each
# cell is explicitly formatted as Text before the new value is written --
# otherwise Excel would happily "helpfully" reinterpret a value such as
# 58575711 as a number (or 03-07-2022 as a date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Row 2 - OrderId
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "58575711"

# Row 3 - OrderId, RecurringExpiry, CloneOrderId
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "58575712"

$ws.Range("AD3").NumberFormat = "@"
$ws.Range("AD3").Value = "03-07-2022"

$ws.Range("R3").NumberFormat = "@"
$ws.Range("R3").Value = "58575713"

# Row 4 - OrderId
$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "58575714"
